$d = $word.ActiveDocument

$oldText = "How can I test to ensure that HeardIT meets the standards of modern applications?"
$part1 = "How can I "
$part2 = "ensure that HeardIT provides its users with sufficient data security mechanisms and follows the modern standards and regulations for data protection"
$part3 = "?"
$newText = $part1 + $part2 + $part3

# Step 1: swap the sentence text in place. Find/Replace keeps the single
# run's existing character formatting (sz=22 / szCs=24), it just rewrites
# the <w:t>.
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# Step 2: locate the paragraph that now carries the replaced sentence so we
# can split its single run into three runs (matching the target XML).
$count = $d.Paragraphs.Count
$target = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "$newText*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $start = $target.Range.Start

    # Toggling a character property (and back) on a sub-range forces Word
    # to break the run at that boundary without altering the visible
    # formatting, which is how the diff ends up with three <w:r> elements
    # that all share the same <w:rPr>.
    $split1 = $d.Range($start, $start + $part1.Length)
    $split1.Bold = 1
    $split1.Bold = 0

    $split2 = $d.Range($start + $part1.Length, $start + $part1.Length + $part2.Length)
    $split2.Bold = 1
    $split2.Bold = 0
}
